$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Row 1429 (year 1979) is unchanged by this edit, but restate its K value explicitly
# (the workbook's stored value has a stray leading space that can otherwise be lost on re-save).
$ws.Range("K1429").Value = 99999999999999991611392.0
$ws.Range("D1400").Value = 10402219140083218432.0
$ws.Range("I1400").Value = -436893203883495168.0
$ws.Range("J1400").Value = -177888112330413.0
$ws.Range("K1400").Value = 104022191400832190119936.0
$ws.Range("D1401").Value = 10388349514563106816.0
$ws.Range("I1401").Value = -436310679611650560.0
$ws.Range("J1401").Value = -177650928180639.0
$ws.Range("K1401").Value = 103883495145631064784896.0
$ws.Range("D1402").Value = 10374479889042995200.0
$ws.Range("I1402").Value = -435728155339805824.0
$ws.Range("J1402").Value = -177413744030866.0
$ws.Range("K1402").Value = 103744798890429956227072.0
$ws.Range("D1403").Value = 10360610263522885632.0
$ws.Range("I1403").Value = -435145631067961216.0
$ws.Range("J1403").Value = -177176559881092.0
$ws.Range("K1403").Value = 103606102635228864446464.0
$ws.Range("D1404").Value = 10346740638002774016.0
$ws.Range("I1404").Value = -434563106796116480.0
$ws.Range("J1404").Value = -176939375731318.0
$ws.Range("K1404").Value = 103467406380027739111424.0
$ws.Range("D1405").Value = 10332871012482662400.0
$ws.Range("I1405").Value = -433980582524271808.0
$ws.Range("J1405").Value = -176702191581544.0
$ws.Range("K1405").Value = 103328710124826630553600.0
$ws.Range("D1406").Value = 10319001386962552832.0
$ws.Range("I1406").Value = -433398058252427264.0
$ws.Range("J1406").Value = -176465007431770.0
$ws.Range("K1406").Value = 103190013869625521995776.0
$ws.Range("D1407").Value = 10305131761442441216.0
$ws.Range("I1407").Value = -432815533980582592.0
$ws.Range("J1407").Value = -176227823281996.0
$ws.Range("K1407").Value = 103051317614424413437952.0
$ws.Range("D1408").Value = 10291262135922329600.0
$ws.Range("I1408").Value = -432233009708737856.0
$ws.Range("J1408").Value = -175990639132222.0
$ws.Range("K1408").Value = 102912621359223288102912.0
$ws.Range("D1409").Value = 10277392510402217984.0
$ws.Range("I1409").Value = -431650485436893184.0
$ws.Range("J1409").Value = -175753454982448.0
$ws.Range("K1409").Value = 102773925104022179545088.0
$ws.Range("D1410").Value = 10263522884882108416.0
$ws.Range("I1410").Value = -431067961165048512.0
$ws.Range("J1410").Value = -175516270832674.0
$ws.Range("K1410").Value = 102635228848821087764480.0
$ws.Range("D1411").Value = 10249653259361998848.0
$ws.Range("I1411").Value = -430485436893203968.0
$ws.Range("J1411").Value = -175279086682901.0
$ws.Range("K1411").Value = 102496532593619995983872.0
$ws.Range("D1412").Value = 10235783633841887232.0
$ws.Range("I1412").Value = -429902912621359296.0
$ws.Range("J1412").Value = -175041902533127.0
$ws.Range("K1412").Value = 102357836338418870648832.0
$ws.Range("D1413").Value = 10221914008321773568.0
$ws.Range("I1413").Value = -429320388349514496.0
$ws.Range("J1413").Value = -174804718383353.0
$ws.Range("K1413").Value = 102219140083217728536576.0
$ws.Range("D1414").Value = 10208044382801664000.0
$ws.Range("I1414").Value = -428737864077669888.0
$ws.Range("J1414").Value = -174567534233579.0
$ws.Range("K1414").Value = 102080443828016636755968.0
$ws.Range("D1415").Value = 10194174757281554432.0
$ws.Range("I1415").Value = -428155339805825280.0
$ws.Range("J1415").Value = -174330350083805.0
$ws.Range("K1415").Value = 101941747572815544975360.0
$ws.Range("D1416").Value = 10180305131761442816.0
$ws.Range("I1416").Value = -427572815533980672.0
$ws.Range("J1416").Value = -174093165934031.0
$ws.Range("K1416").Value = 101803051317614436417536.0
$ws.Range("D1417").Value = 10166435506241331200.0
$ws.Range("I1417").Value = -426990291262135936.0
$ws.Range("J1417").Value = -173855981784257.0
$ws.Range("K1417").Value = 101664355062413311082496.0
$ws.Range("D1418").Value = 10152565880721219584.0
$ws.Range("I1418").Value = -426407766990291264.0
$ws.Range("J1418").Value = -173618797634483.0
$ws.Range("K1418").Value = 101525658807212202524672.0
$ws.Range("D1419").Value = 10138696255201110016.0
$ws.Range("I1419").Value = -425825242718446592.0
$ws.Range("J1419").Value = -173381613484710.0
$ws.Range("K1419").Value = 101386962552011093966848.0
$ws.Range("D1420").Value = 10124826629680998400.0
$ws.Range("I1420").Value = -425242718446601920.0
$ws.Range("J1420").Value = -173144429334936.0
$ws.Range("K1420").Value = 101248266296809985409024.0
$ws.Range("D1421").Value = 10110957004160886784.0
$ws.Range("I1421").Value = -424660194174757184.0
$ws.Range("J1421").Value = -172907245185162.0
$ws.Range("K1421").Value = 101109570041608860073984.0
$ws.Range("D1422").Value = 10097087378640777216.0
$ws.Range("I1422").Value = -424077669902912704.0
$ws.Range("J1422").Value = -172670061035388.0
$ws.Range("K1422").Value = 100970873786407768293376.0
$ws.Range("D1423").Value = 10083217753120665600.0
$ws.Range("I1423").Value = -423495145631067968.0
$ws.Range("J1423").Value = -172432876885614.0
$ws.Range("K1423").Value = 100832177531206659735552.0
$ws.Range("D1424").Value = 10069348127600553984.0
$ws.Range("I1424").Value = -422912621359223296.0
$ws.Range("J1424").Value = -172195692735840.0
$ws.Range("K1424").Value = 100693481276005534400512.0
$ws.Range("D1425").Value = 10055478502080442368.0
$ws.Range("I1425").Value = -422330097087378560.0
$ws.Range("J1425").Value = -171958508586066.0
$ws.Range("K1425").Value = 100554785020804425842688.0
$ws.Range("D1426").Value = 10041608876560332800.0
$ws.Range("I1426").Value = -421747572815534016.0
$ws.Range("J1426").Value = -171721324436292.0
$ws.Range("K1426").Value = 100416088765603334062080.0
$ws.Range("D1427").Value = 10027739251040223232.0
$ws.Range("I1427").Value = -421165048543689408.0
$ws.Range("J1427").Value = -171484140286518.0
$ws.Range("K1427").Value = 100277392510402225504256.0
$ws.Range("D1428").Value = 10013869625520111616.0
$ws.Range("I1428").Value = -420582524271844736.0
$ws.Range("J1428").Value = -171246956136745.0
$ws.Range("K1428").Value = 100138696255201116946432.0
$ws.Range("D1430").Value = 9986130374479888384.0
$ws.Range("I1430").Value = -419417475728155264.0
$ws.Range("J1430").Value = -170772587837197.0
$ws.Range("K1430").Value = 99861303744798883053568.0
$ws.Range("D1431").Value = 9972260748959778816.0
$ws.Range("I1431").Value = -418834951456310784.0
$ws.Range("J1431").Value = -170535403687423.0
$ws.Range("K1431").Value = 99722607489597791272960.0
$ws.Range("D1432").Value = 9958391123439667200.0
$ws.Range("I1432").Value = -418252427184466048.0
$ws.Range("J1432").Value = -170298219537649.0
$ws.Range("K1432").Value = 99583911234396665937920.0
$ws.Range("D1433").Value = 9944521497919555584.0
$ws.Range("I1433").Value = -417669902912621376.0
$ws.Range("J1433").Value = -170061035387875.0
$ws.Range("K1433").Value = 99445214979195557380096.0
$ws.Range("D1434").Value = 9930651872399446016.0
$ws.Range("I1434").Value = -417087378640776704.0
$ws.Range("J1434").Value = -169823851238101.0
$ws.Range("K1434").Value = 99306518723994465599488.0
$ws.Range("D1435").Value = 9916782246879334400.0
$ws.Range("I1435").Value = -416504854368932032.0
$ws.Range("J1435").Value = -169586667088327.0
$ws.Range("K1435").Value = 99167822468793340264448.0
$ws.Range("D1436").Value = 9902912621359222784.0
$ws.Range("I1436").Value = -415922330097087296.0
$ws.Range("J1436").Value = -169349482938553.0
$ws.Range("K1436").Value = 99029126213592231706624.0
$ws.Range("D1437").Value = 9889042995839111168.0
$ws.Range("I1437").Value = -415339805825242688.0
$ws.Range("J1437").Value = -169112298788780.0
$ws.Range("K1437").Value = 98890429958391106371584.0
$ws.Range("D1438").Value = 9875173370319001600.0
$ws.Range("I1438").Value = -414757281553398080.0
$ws.Range("J1438").Value = -168875114639006.0
$ws.Range("K1438").Value = 98751733703190014590976.0
$ws.Range("D1439").Value = 9861303744798892032.0
$ws.Range("I1439").Value = -414174757281553472.0
$ws.Range("J1439").Value = -168637930489232.0
$ws.Range("K1439").Value = 98613037447988922810368.0
$ws.Range("D1440").Value = 9847434119278780416.0
$ws.Range("I1440").Value = -413592233009708736.0
$ws.Range("J1440").Value = -168400746339458.0
$ws.Range("K1440").Value = 98474341192787797475328.0
$ws.Range("D1441").Value = 9833564493758668800.0
$ws.Range("I1441").Value = -413009708737864128.0
$ws.Range("J1441").Value = -168163562189684.0
$ws.Range("K1441").Value = 98335644937586688917504.0
$ws.Range("D1442").Value = 9819694868238557184.0
$ws.Range("I1442").Value = -412427184466019456.0
$ws.Range("J1442").Value = -167926378039910.0
$ws.Range("K1442").Value = 98196948682385563582464.0
$ws.Range("D1443").Value = 9805825242718447616.0
$ws.Range("I1443").Value = -411844660194174848.0
$ws.Range("J1443").Value = -167689193890136.0
$ws.Range("K1443").Value = 98058252427184471801856.0
$ws.Range("D1444").Value = 9791955617198336000.0
$ws.Range("I1444").Value = -411262135922330112.0
$ws.Range("J1444").Value = -167452009740362.0
$ws.Range("K1444").Value = 97919556171983363244032.0
$ws.Range("D1445").Value = 9778085991678224384.0
$ws.Range("I1445").Value = -410679611650485376.0
$ws.Range("J1445").Value = -167214825590589.0
$ws.Range("K1445").Value = 97780859916782237908992.0
$ws.Range("D1446").Value = 9764216366158112768.0
$ws.Range("I1446").Value = -410097087378640704.0
$ws.Range("J1446").Value = -166977641440815.0
$ws.Range("K1446").Value = 97642163661581129351168.0
$ws.Range("D1447").Value = 9750346740638003200.0
$ws.Range("I1447").Value = -409514563106796160.0
$ws.Range("J1447").Value = -166740457291041.0
$ws.Range("K1447").Value = 97503467406380037570560.0
$ws.Range("D1448").Value = 9736477115117891584.0
$ws.Range("I1448").Value = -408932038834951488.0
$ws.Range("J1448").Value = -166503273141267.0
$ws.Range("K1448").Value = 97364771151178912235520.0
$ws.Range("D1449").Value = 9722607489597779968.0
$ws.Range("I1449").Value = -408349514563106752.0
$ws.Range("J1449").Value = -166266088991493.0
$ws.Range("K1449").Value = 97226074895977803677696.0
$ws.Range("D1450").Value = 9708737864077670400.0
$ws.Range("I1450").Value = -407766990291262144.0
$ws.Range("J1450").Value = -166028904841719.0
$ws.Range("K1450").Value = 97087378640776711897088.0
$ws.Range("D1451").Value = 9694868238557558784.0
$ws.Range("I1451").Value = -407184466019417408.0
$ws.Range("J1451").Value = -165791720691945.0
$ws.Range("K1451").Value = 96948682385575586562048.0
$ws.Range("D1452").Value = 9680998613037447168.0
$ws.Range("I1452").Value = -406601941747572800.0
$ws.Range("J1452").Value = -165554536542171.0
$ws.Range("K1452").Value = 96809986130374478004224.0
$ws.Range("D1453").Value = 9667128987517335552.0
$ws.Range("I1453").Value = -406019417475728128.0
$ws.Range("J1453").Value = -165317352392397.0
$ws.Range("K1453").Value = 96671289875173352669184.0
$ws.Range("D1454").Value = 9653259361997225984.0
$ws.Range("I1454").Value = -405436893203883520.0
$ws.Range("J1454").Value = -165080168242624.0
$ws.Range("K1454").Value = 96532593619972260888576.0
$ws.Range("D1455").Value = 9639389736477116416.0
$ws.Range("I1455").Value = -404854368932038848.0
$ws.Range("J1455").Value = -164842984092850.0
$ws.Range("K1455").Value = 96393897364771169107968.0
$ws.Range("D1456").Value = 9625520110957004800.0
$ws.Range("I1456").Value = -404271844660194240.0
$ws.Range("J1456").Value = -164605799943076.0
$ws.Range("K1456").Value = 96255201109570043772928.0
$ws.Range("D1457").Value = 9611650485436893184.0
$ws.Range("I1457").Value = -403689320388349568.0
$ws.Range("J1457").Value = -164368615793302.0
$ws.Range("K1457").Value = 96116504854368935215104.0
$ws.Range("D1458").Value = 9597780859916781568.0
$ws.Range("I1458").Value = -403106796116504832.0
$ws.Range("J1458").Value = -164131431643528.0
$ws.Range("K1458").Value = 95977808599167809880064.0
$ws.Range("D1459").Value = 9583911234396672000.0
$ws.Range("I1459").Value = -402524271844660224.0
$ws.Range("J1459").Value = -163894247493754.0
$ws.Range("K1459").Value = 95839112343966718099456.0
$ws.Range("D1460").Value = 9570041608876560384.0
$ws.Range("I1460").Value = -401941747572815488.0
$ws.Range("J1460").Value = -163657063343980.0
$ws.Range("K1460").Value = 95700416088765609541632.0
$ws.Range("D1461").Value = 9556171983356448768.0
$ws.Range("I1461").Value = -401359223300970816.0
$ws.Range("J1461").Value = -163419879194206.0
$ws.Range("K1461").Value = 95561719833564484206592.0
$ws.Range("D1462").Value = 9542302357836339200.0
$ws.Range("I1462").Value = -400776699029126272.0
$ws.Range("J1462").Value = -163182695044432.0
$ws.Range("K1462").Value = 95423023578363392425984.0
$ws.Range("D1463").Value = 9528432732316227584.0
$ws.Range("I1463").Value = -400194174757281600.0
$ws.Range("J1463").Value = -162945510894659.0
$ws.Range("K1463").Value = 95284327323162283868160.0
$ws.Range("D1464").Value = 9514563106796115968.0
$ws.Range("I1464").Value = -399611650485436864.0
$ws.Range("J1464").Value = -162708326744885.0
$ws.Range("K1464").Value = 95145631067961158533120.0
$ws.Range("D1465").Value = 9500693481276004352.0
$ws.Range("I1465").Value = -399029126213592192.0
$ws.Range("J1465").Value = -162471142595111.0
$ws.Range("K1465").Value = 95006934812760049975296.0
$ws.Range("D1466").Value = 9486823855755894784.0
$ws.Range("I1466").Value = -398446601941747520.0
$ws.Range("J1466").Value = -162233958445337.0
$ws.Range("K1466").Value = 94868238557558941417472.0
$ws.Range("D1467").Value = 9472954230235785216.0
$ws.Range("I1467").Value = -397864077669903040.0
$ws.Range("J1467").Value = -161996774295563.0
$ws.Range("K1467").Value = 94729542302357849636864.0
$ws.Range("D1468").Value = 9459084604715673600.0
$ws.Range("I1468").Value = -397281553398058304.0
$ws.Range("J1468").Value = -161759590145789.0
$ws.Range("K1468").Value = 94590846047156741079040.0
$ws.Range("D1469").Value = 9445214979195561984.0
$ws.Range("I1469").Value = -396699029126213632.0
$ws.Range("J1469").Value = -161522405996015.0
$ws.Range("K1469").Value = 94452149791955615744000.0
$ws.Range("D1470").Value = 9431345353675450368.0
$ws.Range("I1470").Value = -396116504854368896.0
$ws.Range("J1470").Value = -161285221846241.0
$ws.Range("K1470").Value = 94313453536754507186176.0
$ws.Range("D1471").Value = 9417475728155340800.0
$ws.Range("I1471").Value = -395533980582524288.0
$ws.Range("J1471").Value = -161048037696468.0
$ws.Range("K1471").Value = 94174757281553415405568.0
$ws.Range("D1472").Value = 9403606102635229184.0
$ws.Range("I1472").Value = -394951456310679680.0
$ws.Range("J1472").Value = -160810853546694.0
$ws.Range("K1472").Value = 94036061026352290070528.0
$ws.Range("D1473").Value = 9389736477115117568.0
$ws.Range("I1473").Value = -394368932038834944.0
$ws.Range("J1473").Value = -160573669396920.0
$ws.Range("K1473").Value = 93897364771151181512704.0
$ws.Range("D1474").Value = 9375866851595005952.0
$ws.Range("I1474").Value = -393786407766990272.0
$ws.Range("J1474").Value = -160336485247146.0
$ws.Range("K1474").Value = 93758668515950056177664.0
$ws.Range("D1475").Value = 9361997226074896384.0
$ws.Range("I1475").Value = -393203883495145600.0
$ws.Range("J1475").Value = -160099301097372.0
$ws.Range("K1475").Value = 93619972260748964397056.0
